# Scheduled runner update: refresh market-board derived price/profit figures
# across several worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) of the
# Pandaemonium profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 641.8570999999999
$ws.Range("I12").Value = 270.5
$ws.Range("J12").Value = 790.4
$ws.Range("K12").Value = 270.5
$ws.Range("L12").Value = 790.4
$ws.Range("M12").Value = -100.5
$ws.Range("N12").Value = -1130.4

$ws.Range("H70").Value = 1714.6578
$ws.Range("I70").Value = 1644.6111
$ws.Range("J70").Value = 1777.7
$ws.Range("K70").Value = 4933.8333
$ws.Range("L70").Value = 5333.1
$ws.Range("M70").Value = -4663.8333
$ws.Range("N70").Value = -5873.1

$ws.Range("H73").Value = 1714.6578
$ws.Range("I73").Value = 1644.6111
$ws.Range("J73").Value = 1777.7
$ws.Range("K73").Value = 4933.8333
$ws.Range("L73").Value = 5333.1
$ws.Range("M73").Value = -3997.8333
$ws.Range("N73").Value = -7205.1

$ws.Range("H88").Value = 1862.5
$ws.Range("I88").Value = 816.6667
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 816.6667
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -410.6667
$ws.Range("N88").Value = -5812

$ws.Range("H91").Value = 1862.5
$ws.Range("I91").Value = 816.6667
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 816.6667
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = 587.3333
$ws.Range("N91").Value = -7808

$ws.Range("H135").Value = 65218596
$ws.Range("I135").Value = 29412326
$ws.Range("J135").Value = 166669700
$ws.Range("K135").Value = 264710934
$ws.Range("L135").Value = 1500027300
$ws.Range("M135").Value = -264708399
$ws.Range("N135").Value = -1500032370

$ws.Range("H137").Value = 760485.25
$ws.Range("I137").Value = 3080.9
$ws.Range("J137").Value = 1391655.5
$ws.Range("K137").Value = 9242.700000000001
$ws.Range("L137").Value = 4174966.5
$ws.Range("M137").Value = -6692.700000000001
$ws.Range("N137").Value = -4180066.5

$ws.Range("H138").Value = 3219.5
$ws.Range("I138").Value = 3248.125
$ws.Range("J138").Value = 3214.92
$ws.Range("K138").Value = 9744.375
$ws.Range("L138").Value = 9644.76
$ws.Range("M138").Value = -4604.375
$ws.Range("N138").Value = -19924.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 77145
$ws.Range("J121").Value = 77145
$ws.Range("L121").Value = 77145
$ws.Range("N121").Value = -80639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 309.2857
$ws.Range("I22").Value = 325.45456
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 325.45456
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 24.54543999999999
$ws.Range("N22").Value = -950

$ws.Range("H86").Value = 2259.7
$ws.Range("I86").Value = 2033.3334
$ws.Range("J86").Value = 2599.25
$ws.Range("K86").Value = 2033.3334
$ws.Range("L86").Value = 2599.25
$ws.Range("M86").Value = -910.3334
$ws.Range("N86").Value = -4845.25

$ws.Range("H89").Value = 2259.7
$ws.Range("I89").Value = 2033.3334
$ws.Range("J89").Value = 2599.25
$ws.Range("K89").Value = 10166.667
$ws.Range("L89").Value = 12996.25
$ws.Range("M89").Value = -4550.666999999999
$ws.Range("N89").Value = -24228.25

$ws.Range("H132").Value = 2351.2793
$ws.Range("I132").Value = 1649.88
$ws.Range("K132").Value = 4949.64
$ws.Range("M132").Value = -2419.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6671047.5
$ws.Range("I5").Value = 384.6316
$ws.Range("J5").Value = 27794814
$ws.Range("K5").Value = 1153.8948
$ws.Range("L5").Value = 83384442
$ws.Range("M5").Value = -1041.8948
$ws.Range("N5").Value = -83384666

$ws.Range("H68").Value = 3029.419
$ws.Range("J68").Value = 4665.6855
$ws.Range("L68").Value = 13997.0565
$ws.Range("N68").Value = -15619.0565

$ws.Range("H71").Value = 3029.419
$ws.Range("J71").Value = 4665.6855
$ws.Range("L71").Value = 41991.1695
$ws.Range("N71").Value = -50103.1695

$ws.Range("H125").Value = 1747.4
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H135").Value = 6671047.5
$ws.Range("I135").Value = 384.6316
$ws.Range("J135").Value = 27794814
$ws.Range("K135").Value = 3461.6844
$ws.Range("L135").Value = 250153326
$ws.Range("M135").Value = -926.6844000000001
$ws.Range("N135").Value = -250158396

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12999.875
$ws.Range("I5").Value = 8666.666999999999
$ws.Range("J5").Value = 15599.8
$ws.Range("K5").Value = 8666.666999999999
$ws.Range("L5").Value = 15599.8
$ws.Range("M5").Value = -8554.666999999999
$ws.Range("N5").Value = -15823.8

$ws.Range("H123").Value = 10334.125
$ws.Range("J123").Value = 10334.125
$ws.Range("L123").Value = 10334.125
$ws.Range("N123").Value = -15234.125

$ws.Range("H132").Value = 9299
$ws.Range("I132").Value = 1805.1111
$ws.Range("J132").Value = 22788
$ws.Range("K132").Value = 5415.3333
$ws.Range("L132").Value = 68364
$ws.Range("M132").Value = -2885.3333
$ws.Range("N132").Value = -73424

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3000000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393

$ws.Range("H68").Value = 4017.5715
$ws.Range("I68").Value = 2466.5
$ws.Range("K68").Value = 2466.5
$ws.Range("M68").Value = -1717.5

$ws.Range("H71").Value = 4017.5715
$ws.Range("I71").Value = 2466.5
$ws.Range("K71").Value = 12332.5
$ws.Range("M71").Value = -8588.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("J5").Value = 10000000
$ws.Range("L5").Value = 10000000
$ws.Range("N5").Value = -10000224

$ws.Range("H13").Value = 3131.5
$ws.Range("I13").Value = 3131.5
$ws.Range("K13").Value = 3131.5
$ws.Range("M13").Value = -2991.5

$ws.Range("H23").Value = 14124.286
$ws.Range("I23").Value = 3467.5
$ws.Range("J23").Value = 28333.334
$ws.Range("K23").Value = 3467.5
$ws.Range("L23").Value = 28333.334
$ws.Range("M23").Value = -3238.5
$ws.Range("N23").Value = -28791.334

$ws.Range("H123").Value = 38378.875
$ws.Range("J123").Value = 38378.875
$ws.Range("L123").Value = 38378.875
$ws.Range("N123").Value = -48178.875

$ws.Range("H132").Value = 1554.3829
$ws.Range("I132").Value = 1815.2963
$ws.Range("J132").Value = 1202.15
$ws.Range("K132").Value = 5445.8889
$ws.Range("L132").Value = 3606.45
$ws.Range("M132").Value = -2915.8889
$ws.Range("N132").Value = -8666.450000000001
